# Update the experimental-data-density sheet: the ammonia partial
# pressure (G), mass-transfer (H) and diffusion-flux (I) formulas were
# reworked.  All three columns use shared formulas in three separate
# blocks (rows 2, 3:66, 67:130, 131:134) so we re-apply the formula to
# the first cell of each block (row 2) and to each shared-formula
# anchor cell (rows 3, 67, 131) with its own range, which keeps the
# existing shared-formula grouping intact while updating every
# dependent cell's calculated value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G: 0.01*101325/D{row}*100  ->  0.015*101325/D{row}*100
$ws.Range("G2").Formula = "=0.015*101325/D2*100"
$ws.Range("G3:G66").Formula = "=0.015*101325/D3*100"
$ws.Range("G67:G130").Formula = "=0.015*101325/D67*100"
$ws.Range("G131:G134").Formula = "=0.015*101325/D131*100"

# Column H: 0.001/E{row}*100  ->  0.001/B{row}*100
$ws.Range("H2").Formula = "=0.001/B2*100"
$ws.Range("H3:H66").Formula = "=0.001/B3*100"
$ws.Range("H67:H130").Formula = "=0.001/B67*100"
$ws.Range("H131:H134").Formula = "=0.001/B131*100"

# Column I: 0.5*0.000000001*F{row}/(9*0.001)  ->
#           (0.2/9000+F{row}/(9*0.001)*0.5*0.000000001)*100
$ws.Range("I2").Formula = "=(0.2/9000+F2/(9*0.001)*0.5*0.000000001)*100"
$ws.Range("I3:I66").Formula = "=(0.2/9000+F3/(9*0.001)*0.5*0.000000001)*100"
$ws.Range("I67:I130").Formula = "=(0.2/9000+F67/(9*0.001)*0.5*0.000000001)*100"
$ws.Range("I131:I134").Formula = "=(0.2/9000+F131/(9*0.001)*0.5*0.000000001)*100"

# The active selection moved from I2 to G2.
$ws.Range("G2").Select()
